# Daily attendance processing - 2026-01-07 08:42:53
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Class Statistics summary (K/L column) ---------------------------------
# Missing Sessions: 21 -> 27
$ws.Range("L7").Value = 27
# Pending Sessions: 84 -> 78
$ws.Range("L8").Value = 78

# --- "Recorded By" cells: System moved to front of the list ----------------
$recordedByRows = @(8,9,10,12,14,15,17,34,35,36,38,40,41,43,60,61,62,64,66,67,69,86,87,88,90,92,93,95,112,113,114,116,118,119,121,138,139,140,142,144,145,147)
foreach ($r in $recordedByRows) {
    $ws.Range("G$r").Value = "System, dnasr281@gmail.com"
}

# --- Group statistics table (rows 21-26): one session per group moved from
#     "Pending" to "Missing" (Missing +1, Pending -1) ------------------------
$groupStatRows = @(21,22,23,24,25,26)
foreach ($r in $groupStatRows) {
    $ws.Range("P$r").Value = $ws.Range("P$r").Value() + 1
    $ws.Range("Q$r").Value = $ws.Range("Q$r").Value() - 1
}

# --- Session rows that flipped from "Pending" (future/yellow) to
#     "Not Recorded" (past-due/pink) now that 07/01/2026 has elapsed ---------
# Each target row already has a same-shaped sibling row formatted as
# "Not Recorded" (the row directly above it) - copy that formatting across,
# then restore the correct session status text.
$pendingToMissingRows = @(177,204,231,258,285,312)
foreach ($r in $pendingToMissingRows) {
    $srcRow = $r - 1
    $ws.Range("A${srcRow}:I${srcRow}").Copy()
    $ws.Range("A${r}:I${r}").PasteSpecial(-4122)
    $ws.Range("I$r").Value = "Not Recorded"
}
$excel.CutCopyMode = $false
